$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.074.80"
$ws.Range("E2").Value = "  +0.27%  "

$ws.Range("D3").Value = "1.833.26"
$ws.Range("E3").Value = "  +0.28%  "

$ws.Range("D4").Value = "'0.9976"
$ws.Range("E4").Value = "  +0.17%  "

$ws.Range("D5").Value = "'242.38"
$ws.Range("E5").Value = "  -0.60%  "

$ws.Range("E6").Value = "  -2.01%  "

$ws.Range("D7").Value = "'1.001"
$ws.Range("E7").Value = "  +0.27%  "

$ws.Range("D8").Value = "'0.07447"
$ws.Range("E8").Value = "  -0.72%  "

$ws.Range("D9").Value = "'0.2930"
$ws.Range("E9").Value = "  -0.40%  "

$ws.Range("D10").Value = "'23.06"
$ws.Range("E10").Value = "  -0.54%  "

$ws.Range("D11").Value = "'0.07662"
$ws.Range("E11").Value = "  -0.52%  "

$ws.Range("D12").Value = "1.851.00"
$ws.Range("E12").Value = "  +1.20%  "

$ws.Range("D13").Value = "'5.012"
$ws.Range("E13").Value = "  +0.39%  "

$ws.Range("E14").Value = "  +0.92%  "

$ws.Range("D15").Value = "'82.93"
$ws.Range("E15").Value = "  -0.24%  "

$ws.Range("D16").Value = "'0.000009151"
$ws.Range("E16").Value = "  -6.39%  "

$ws.Range("D17").Value = "'5.900"
$ws.Range("E17").Value = "  -2.23%  "

$ws.Range("D18").Value = "29.097.88"
$ws.Range("E18").Value = "  +0.26%  "

$ws.Range("D19").Value = "2.089.32"
$ws.Range("E19").Value = "  +0.89%  "

$ws.Range("D20").Value = "'240.79"
$ws.Range("E20").Value = "  +6.64%  "

$ws.Range("D21").Value = "'12.71"
$ws.Range("E21").Value = "  +1.07%  "

$ws.Range("D22").Value = "'1.002"
$ws.Range("E22").Value = "  +0.34%  "

$ws.Range("D23").Value = "'7.215"
$ws.Range("E23").Value = "  +1.14%  "

$ws.Range("D24").Value = "'0.9991"
$ws.Range("E24").Value = "  +0.12%  "

$ws.Range("D25").Value = "'159.06"
$ws.Range("E25").Value = "  -1.00%  "

$ws.Range("D26").Value = "'0.1414"
$ws.Range("E26").Value = "  -0.01%  "

$ws.Range("D27").Value = "'8.511"
$ws.Range("E27").Value = "  -0.07%  "

$ws.Range("D28").Value = "'17.89"
$ws.Range("E28").Value = "  -0.25%  "

$ws.Range("D29").Value = "'1.500"
$ws.Range("E29").Value = "  +0.04%  "

$ws.Range("D30").Value = "'0.05590"
$ws.Range("E30").Value = "  +2.10%  "

$ws.Range("E31").Value = "  +1.62%  "

$ws.Range("D32").Value = "'4.131"
$ws.Range("E32").Value = "  +0.06%  "

$ws.Range("E33").Value = "  +0.91%  "

$ws.Range("E34").Value = "  -0.85%  "

$ws.Range("D35").Value = "'0.7428"
$ws.Range("E35").Value = "  -0.24%  "

$ws.Range("E36").Value = "  +0.51%  "

$ws.Range("D37").Value = "'2.655"
$ws.Range("E37").Value = "  +1.76%  "

$ws.Range("D38").Value = "'2.770"
$ws.Range("E38").Value = "  +0.75%  "

$ws.Range("D39").Value = "'0.01785"
$ws.Range("E39").Value = "  +0.06%  "

$ws.Range("D40").Value = "1.210.96"
$ws.Range("E40").Value = "  -2.44%  "

$ws.Range("D41").Value = "'6.404"
$ws.Range("E41").Value = "  -4.58%  "

$ws.Range("D42").Value = "'0.8967"
$ws.Range("E42").Value = "  -0.67%  "

$ws.Range("E43").Value = "  +0.21%  "

$ws.Range("D44").Value = "'101.46"
$ws.Range("E44").Value = "  -0.20%  "

$ws.Range("D45").Value = "1.984.52"
$ws.Range("E45").Value = "  +0.71%  "

$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").Value = "'65.49"
$ws.Range("E46").Value = "  +0.43%  "

$ws.Range("B47").Value = "BabyDogeCoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D47").Value = "'0.00000000122"
$ws.Range("E47").Value = "  -2.46%  "

$ws.Range("D48").Value = "'0.5086"
$ws.Range("E48").Value = "  +0.29%  "

$ws.Range("D49").Value = "'0.4065"
$ws.Range("E49").Value = "  +0.29%  "

$ws.Range("D50").Value = "'9.144"
$ws.Range("E50").Value = "  +2.03%  "

$ws.Range("D51").Value = "'0.05841"
$ws.Range("E51").Value = "  +0.91%  "
